# Fill in the requirements table (rows 2-5) with the RF06..RF09 entries and
# their descriptions, matching the order the author originally typed them in
# (codes for rows 3-5 first, then the descriptions out of row order) so the
# shared-string table comes out in the same sequence as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "RF06"
$ws.Range("B2").Value = "Se requiere que el usuario confirme el cierre del pedido"

$ws.Range("A3").Value = "RF07"
$ws.Range("A4").Value = "RF08"
$ws.Range("A5").Value = "RF09"

$ws.Range("B3").Value = "El usuario tiene que estar logueado para efectivizar el pedido"
$ws.Range("B5").Value = "Se requiere notificación al usuario del pedido realizado"
$ws.Range("B4").Value = "Se requiere que el usuario pueda modificar la dirección en la cuál se enviará el pedido"

# B5 ends up with its own (slightly different/"automatic" colored) font.
$ws.Range("B5").Font.ThemeColor = 1

# Row heights settle to the auto-fit values once the wrapped text is in place.
$ws.Rows(2).RowHeight = 12.75
$ws.Rows(3).RowHeight = 16.5
$ws.Rows(4).RowHeight = 33
$ws.Rows(5).RowHeight = 13.5

# The former blank row 6 is removed entirely (dimension shrinks to A1:B5).
$ws.Rows(6).Delete()

# Selection moves on to the next empty row in column B.
$ws.Range("B7").Select() | Out-Null
